# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the zh-cn and de-de
# sheets to reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2:E3").Value = "2016-03-15 04:06:54"
$zhcn.Range("H2:H3").Value = "2016-03-15 04:07:38"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2:E3").Value = "2016-03-15 04:07:01"
$dede.Range("H2:H3").Value = "2016-03-15 04:07:51"
